$d = $word.ActiveDocument

# NOTE: edit 3 is applied first, while the trailing ";<newline>" at the end
# of the paragraph is still the *only* occurrence of that two-character
# sequence in the document. (Edit 1 below introduces a second ";<newline>"
# sequence earlier in the same paragraph, which would otherwise shadow the
# Find() match for edit 3.)

# ------------------------------------------------------------------
# Edit 3: ";\n" -> ";"
#   The trailing line break at the end of the paragraph is removed.
# ------------------------------------------------------------------
$r3 = $d.Content
$r3.Find.Execute(";`n") | Out-Null
$charStart3 = $r3.End - 1
$charEnd3 = $charStart3 + 1
$rr3 = $d.Range($charStart3, $charEnd3)
$rr3.Delete()

# ------------------------------------------------------------------
# Edit 1: "; replaced: " -> ";\nreplaced: "
#   The space right after the semicolon becomes a line break.
# ------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("; replaced: ") | Out-Null
$charStart1 = $r1.Start + 1
$charEnd1 = $charStart1 + 1
$rr1 = $d.Range($charStart1, $charEnd1)
$rr1.Text = "`n"

# ------------------------------------------------------------------
# Edit 2: "; noted: " -> "; noted:\n"
#   The trailing space becomes a line break.
# ------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("; noted: ") | Out-Null
$charStart2 = $r2.End - 1
$charEnd2 = $charStart2 + 1
$rr2 = $d.Range($charStart2, $charEnd2)
$rr2.Text = "`n"

Write-Output "edits applied"
